# Update the final paragraph of the DML section (language tag fix) and
# append the newly authored PL/SQL content described in the commit
# "updated documentatie - finished bazele teoretice".

$d = $word.ActiveDocument

# Locate the paragraph that currently ends the document body by searching
# for its distinctive trailing sentence, then resolve it to the owning
# Paragraph so we get the whole <w:p> (including its paragraph mark).
$locator = $d.Content
$locator.Find.ClearFormatting()
$locator.Find.Text = "grupate după diferite criterii, sau prin diferite metode de selecție."
$found = $locator.Find.Execute()
if (-not $found) {
    throw "Could not locate the target paragraph."
}
$targetParagraph = $locator.Paragraphs(1)
$targetRange = $targetParagraph.Range

# Remove the whole paragraph (its text is re-supplied, with a corrected
# paragraph-mark language, by the XML block inserted below) and collapse
# the insertion point to where it used to start.
$insertStart = $targetRange.Start
$targetRange.Delete()
$insertPoint = $d.Range($insertStart, $insertStart)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="6A2BCCEC" w14:textId="7912ADAC" w:rsidR="00B516F9" w:rsidRPr="00C66B38" w:rsidRDefault="00B516F9" w:rsidP="00254936"><w:pPr><w:rPr><w:lang w:val="ro-RO"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:tab/><w:t>Din cadrul instrucțiunilor DML ne vom folosi de INSERT, UPDATE, DELETE, SELECT și sub interogări. Cu ajutorul instrucțiunii INSERT vom adaugă date, cu UPDATE le vom modifica</w:t></w:r><w:r w:rsidR="00615226"><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t xml:space="preserve"> iar</w:t></w:r><w:r><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00615226"><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t>prin DELETE le ștergem. Cu ajutorul sub interogărilor putem să vizualizăm datele în diferite forme (ascendent/descendent</w:t></w:r><w:r w:rsidR="00C66B38"><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t>),</w:t></w:r><w:r w:rsidR="00615226"><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00C66B38"><w:rPr><w:lang w:val="ro-RO"/></w:rPr><w:t>grupate după diferite criterii, sau prin diferite metode de selecție.</w:t></w:r></w:p><w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="ro-RO"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="ro-RO"/>
        </w:rPr>
        <w:tab/>
        <w:t>Putem de asemenea să vedem</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="ro-RO"/>
        </w:rPr>
        <w:t xml:space="preserve"> datele </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="ro-RO"/>
        </w:rPr>
        <w:t>comune din două tabele</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="ro-RO"/>
        </w:rPr>
        <w:t xml:space="preserve"> (</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="ro-RO"/>
        </w:rPr>
        <w:t>intersecția</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="ro-RO"/>
        </w:rPr>
        <w:t xml:space="preserve">) </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="ro-RO"/>
        </w:rPr>
        <w:t xml:space="preserve">sau diferențele între ele </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="ro-RO"/>
        </w:rPr>
        <w:t>(</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="ro-RO"/>
        </w:rPr>
        <w:t>diferența</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="ro-RO"/>
        </w:rPr>
        <w:t>) folosind JOIN-uri. INNER se folosește pentru prima situație iar LEFT/RIGHT pentru cea de-a doua.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="ro-RO"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="ro-RO"/>
        </w:rPr>
        <w:tab/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="ro-RO"/>
        </w:rPr>
        <w:t>PL/SQL este o extensie a limbajului SQL care permite pe lângă proprietățile limbajului SQL declararea de variabile și constante, controlul fluxului, declararea de proceduri și funcții și multe altele.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="ro-RO"/>
        </w:rPr>
        <w:t xml:space="preserve">  Pentru a putea scrie proceduri, funcții și altele de </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="ro-RO"/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t>care avem nevoie, trebuie să ținem cont de mai multe lucruri</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">: de structura unui bloc anonim, de tipuri de date </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="ro-RO"/>
        </w:rPr>
        <w:t>în PL/SQL, structuri de decizie și iterare. Blocurile anonime în PL/SQL au un DECLARE care poate fi opțional, BEGIN care e obligatoriu – sub acesta fiind comenzi SQL și instrucțiuni PL/SQL</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>; EXCEPTION, op</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="ro-RO"/>
        </w:rPr>
        <w:t>țional pentru acțiuni executate în caz de ridicare de excepții și END care este obligatoriu.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="ro-RO"/>
        </w:rPr>
        <w:tab/>
        <w:t>Există mai multe tipuri de date în PL/SQL și anume</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>: scalare (NUMBER, CHARACTER, DATE, BOOLEAN)</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>, compuse (</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="ro-RO"/>
        </w:rPr>
        <w:t>înregistrări), referință (</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>de exemplu REF CURSOR) sau obiecte mari – sunt niște indicatori către obiecte mari stocate separate de alte date (imagini grafice, text, clipuri video)</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:tab/>
        <w:t>Structurile de decizie presupun specificarea uneia sau mai multor condiții unor instrucțiuni care urmează să fie executate dacă evaluarea condiției are valoarea adevărat</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>. Opțional se pot defini instrucțiuni care se execută în cazul valorii de fals a condiției.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> În PL/SQL sunt disponibile: IF-THEN, IF-THEN-ELSE (de asemenea </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="ro-RO"/>
        </w:rPr>
        <w:t>și imbricate)</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>, IF-THEN-ELSIF, CASE, CASE (searched).</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="ro-RO"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:tab/>
        <w:t>Structurile de iterare presupun execu</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="ro-RO"/>
        </w:rPr>
        <w:t>ția de mai multe ori a instrucțiunilor definite într-un bloc de cod. În PL/SQL sunt disponibile LOOP, WHILE, FOR sau imbricări între acestea.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="ro-RO"/>
        </w:rPr>
        <w:t xml:space="preserve"> În ajutorul acestor structuri avem și cele de control ale iterării</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>: EXIT (se situeaz</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="ro-RO"/>
        </w:rPr>
        <w:t>ă în buclă și execuția se termină imediat, controlul programului se reia cu prima instrucțiune de după buclă), EXIT WHEN (se iese dacă condiția se evaluează ca true), CONTINUE (forțează ca următoarea iterație să aibă loc, actuala oprindu-se la întâlnirea cuvântului CONTINUE)</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="ro-RO"/>
        </w:rPr>
        <w:t xml:space="preserve"> și GOTO (oferă un salt necondiționat la o instrucțiune etichetată din </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="ro-RO"/>
        </w:rPr>
        <w:t xml:space="preserve">același </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="ro-RO"/>
        </w:rPr>
        <w:t>subprogram)</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="ro-RO"/>
        </w:rPr>
        <w:t xml:space="preserve">. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="ro-RO"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="ro-RO"/>
        </w:rPr>
        <w:tab/>
        <w:t>Subprogramele pot fi de două tipuri – proceduri care se utilizează pentru a efectua o acțiune, funcții care se utilizează pentru a calcula și returna o valoare</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>;</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="ro-RO"/>
        </w:rPr>
        <w:t xml:space="preserve"> și pot fi locale – definite în cadrul altui bloc PL/SQL sau subprogram sau stocate  - create folosind comanda CREATE.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="ro-RO"/>
        </w:rPr>
        <w:t xml:space="preserve"> Asemeni blocurilor anonime au trei părți</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>: una declarativ</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="ro-RO"/>
        </w:rPr>
        <w:t>ă care are declarații de variabile, tipuri, constante</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>; una executabil</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="ro-RO"/>
        </w:rPr>
        <w:t>ă – obligatorie și care conține instrucțiunile care efectuează acțiunea dorită</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>; tratarea excep</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="ro-RO"/>
        </w:rPr>
        <w:t xml:space="preserve">țiilor </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>– care con</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="ro-RO"/>
        </w:rPr>
        <w:t>ține codul gestionării erorilor</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="ro-RO"/>
        </w:rPr>
        <w:t xml:space="preserve"> de rulare</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="ro-RO"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="ro-RO"/>
        </w:rPr>
        <w:t xml:space="preserve"> Procedurile stocate devin obiecte în schema utilizatorului care a creat-o putând fi apelate în mod explicit de către un client SQL sau dintr-o altă secvență de cod.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="ro-RO"/>
        </w:rPr>
        <w:t xml:space="preserve"> Funcțiile stocate sunt utilizate pentru a calcula și returna o valoare iar corpul acestora trebuie să conțină cel puțin o comandă RETURN a unei date având tipul specificat)</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="ro-RO"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$insertPoint.InsertXML($xml)

Write-Output ("Paragraphs after edit: " + $d.Paragraphs.Count)
